# 17.1.1 - add a 2023 data column (T) to the revenue table and refresh
# the 2022 column's preliminary figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- column widths: D:L, M:O, P:T get an explicit (narrower) width ---
$ws.Range("D1:L1").ColumnWidth = 8.5703125
$ws.Range("M1:O1").ColumnWidth = 8.5703125
$ws.Range("P1:T1").ColumnWidth = 8.5703125

# --- header row is a touch taller now that it wraps over two columns ---
$ws.Rows.Item(1).RowHeight = 42.75

# --- revise the (previously preliminary) 2022 column values ... ---
$ws.Cells.Item(5, 19).Value = 29.5
$ws.Cells.Item(6, 19).Value = 22.4
$ws.Cells.Item(8, 19).Value = 1.9
$ws.Cells.Item(9, 19).Value = 5.0999999999999996

# --- ... and populate the new 2023 column (T) ---
# Copy each row's 2022-column (S) formatting into the new T cell first so
# the new column inherits the same number formats / fonts / borders.
$rows = 4,5,6,7,8,9,10
foreach ($r in $rows) {
    $ws.Cells.Item($r, 19).Copy()
    $ws.Cells.Item($r, 20).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Cells.Item(4, 20).Value = 2023
$ws.Cells.Item(5, 20).Value = 29.4
$ws.Cells.Item(6, 20).Value = 22.1
$ws.Cells.Item(7, 20).Value = "-"
$ws.Cells.Item(8, 20).Value = 1.2
$ws.Cells.Item(9, 20).Value = 6.1
$ws.Cells.Item(10, 20).Value = 0

# --- reset the saved selection back to the top-left cell ---
$ws.Range("A1").Select()
